$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1698171025303532
$ws.Range("C2").Value = 0.1698171025303532
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 20

# Row 3
$ws.Range("B3").Value = 0.5424796504655186
$ws.Range("C3").Value = 0.5424796504655186

# Row 4
$ws.Range("B4").Value = 0.7170038108021736
$ws.Range("C4").Value = 0.7170038108021736

# Row 5
$ws.Range("B5").Value = 0.009443920072667884
$ws.Range("C5").Value = 0.005099293730101737
$ws.Range("D5").Value = 231
$ws.Range("E5").Value = 192

# Row 6
$ws.Range("B6").Value = 0.9391719980261437
$ws.Range("C6").Value = 0.9391719980261437

# Row 7
$ws.Range("B7").Value = 0.00005233867296621407
$ws.Range("C7").Value = 0.00003642386899911322
$ws.Range("D7").Value = 865
$ws.Range("E7").Value = 800

# Row 8
$ws.Range("B8").Value = 0.1001174851394663
$ws.Range("C8").Value = 0.08855716142974797
$ws.Range("D8").Value = 136
$ws.Range("E8").Value = 124

# Row 9
$ws.Range("B9").Value = 0.03547247095417511
$ws.Range("C9").Value = 0.031028515146612
$ws.Range("D9").Value = 293
$ws.Range("E9").Value = 269

# Row 10
$ws.Range("B10").Value = 0.3880699071844666
$ws.Range("C10").Value = 0.3763010113054139
$ws.Range("D10").Value = 42
$ws.Range("G10").Value = 127

# Row 11
$ws.Range("B11").Value = 0.000000003872273723693762
$ws.Range("C11").Value = 0.000000002136035664968034
$ws.Range("D11").Value = 1116
$ws.Range("E11").Value = 1056
$ws.Range("F11").Value = 1484
$ws.Range("G11").Value = 1480

# Row 12
$ws.Range("B12").Value = 0.07952597316640113
$ws.Range("C12").Value = 0.07952597316640113

# Row 13
$ws.Range("B13").Value = 0.5035560325870317
$ws.Range("C13").Value = 0.5035560325870316

# Row 14
$ws.Range("B14").Value = 0.7688953153587761
$ws.Range("C14").Value = 0.7617814239549737
$ws.Range("E14").Value = 25

# Row 15
$ws.Range("B15").Value = 0.7894068617535623
$ws.Range("C15").Value = 0.7868686898597245
$ws.Range("D15").Value = 18
$ws.Range("G15").Value = 31

# Row 16
$ws.Range("B16").Value = 0.6266129947137067
$ws.Range("C16").Value = 0.6266129947137067

# Row 17
$ws.Range("B17").Value = 0.6448527714974646
$ws.Range("C17").Value = 0.6448527714974646
